$wb = $excel.ActiveWorkbook

# --- ODI Batting: the (blank) INNING_NUMBER cells in B2/B3 carry no data;
#     clear them so they are dropped entirely, matching the source export. ---
$wsBatting = $wb.Worksheets.Item("ODI Batting")
$wsBatting.Range("B2").ClearContents()
$wsBatting.Range("B3").ClearContents()

# --- Add the new "ODI Batting Extra" sheet at the end of the workbook ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsExtra = $wb.Worksheets.Add($null, $lastSheet)
$wsExtra.Name = "ODI Batting Extra"

# Header row (text)
$wsExtra.Range("A1").Value = "MATCH_CODE"
$wsExtra.Range("B1").Value = "BATTING_POSITION"
$wsExtra.Range("C1").Value = "NUM_4"
$wsExtra.Range("D1").Value = "NUM_6"
$wsExtra.Range("E1").Value = "PERCENT_RUNS_OF_TOTAL"
$wsExtra.Range("F1").Value = "MAN_OF_MATCH"

# Reuse the bold/bordered header style already used on the other sheets
$wsPlayerInfo = $wb.Worksheets.Item("Player Info")
$wsPlayerInfo.Range("A1").Copy()
$wsExtra.Range("A1:F1").PasteSpecial(-4122)  # xlPasteFormats

# Data rows. MATCH_CODE is a text field (like on "ODI Batting"/"ODI Bowling"),
# so copy it across as a value from an existing text cell to avoid Excel's
# automatic text->number coercion on a literal numeric-looking string.
$wsBatting.Range("D2").Copy()  # "4293" stored as text
$wsExtra.Range("A2").PasteSpecial(-4163)  # xlPasteValues

$wsBatting.Range("D3").Copy()  # "4295" stored as text
$wsExtra.Range("A3").PasteSpecial(-4163)  # xlPasteValues

$wsExtra.Range("B2").Value = 8

$wsExtra.Range("F2").Value = "NO"
$wsExtra.Range("F3").Value = "NO"

$excel.CutCopyMode = $false
